{"js": "// Remove the \"Introducci\u00f3n\" section's illustrative screenshot (a 404 page\n// capture, rId10 / image2.png) and the paragraph of body text that\n// immediately followed it, while keeping the \"Introducci\u00f3n\" heading and the\n// subsequent \"El resultado es un sistema...\" paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Introducci\u00f3n\" heading paragraph.\nconst introIndex = paragraphs.items.findIndex(\n  (p) => p.text.trim() === \"Introducci\u00f3n\"\n);\n\nif (introIndex === -1 || introIndex + 2 >= paragraphs.items.length) {\n  throw new Error('Could not locate the \"Introducci\u00f3n\" paragraph and its two following paragraphs.');\n}\n\n// The paragraph right after the heading holds the inline picture, and the\n// one after that holds the \"Este proyecto explora...\" text. Delete both,\n// leaving the heading and the paragraph that originally came after them\n// (\"El resultado es un sistema...\") intact.\nconst imageParagraph = paragraphs.items[introIndex + 1];\nconst textParagraph = paragraphs.items[introIndex + 2];\n\ntextParagraph.delete();\nimageParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the \"Introducci\u00f3n\" section's illustrative screenshot (a 404 page\n# capture, rId10 / image2.png) and the paragraph of body text that\n# immediately followed it, while keeping the \"Introducci\u00f3n\" heading and the\n# subsequent \"El resultado es un sistema...\" paragraph untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$introIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Introducci\u00f3n\") {\n        $introIndex = $i\n        break\n    }\n}\n\nif ($introIndex -eq -1 -or ($introIndex + 2) -gt $count) {\n    throw \"Could not locate the 'Introducci\u00f3n' paragraph and its two following paragraphs.\"\n}\n\n# Paragraph right after the heading holds the inline picture; the one after\n# that holds the \"Este proyecto explora...\" text. Delete the later one first\n# so the earlier index stays valid.\n$d.Paragraphs.Item($introIndex + 2).Range.Delete()\n$d.Paragraphs.Item($introIndex + 1).Range.Delete()\n"}
